$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 420, shifting existing rows 420-457 down to 421-458
$ws.Rows.Item(420).Insert()

# Populate the newly inserted row 420 with data (new weekly record)
$ws.Cells.Item(420, 1).Value = 4
$ws.Cells.Item(420, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(420, 3).Value = "Los Lagos"
$ws.Cells.Item(420, 4).Value = 45106
$ws.Cells.Item(420, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(420, 5).Value = 10
$ws.Cells.Item(420, 6).Value = 100112043
$ws.Cells.Item(420, 7).Value = "Pepino ensalada"
$ws.Cells.Item(420, 8).Value = "Sin especificar"
$ws.Cells.Item(420, 9).Value = "Primera"
$ws.Cells.Item(420, 10).Value = 120
$ws.Cells.Item(420, 11).Value = 18000
$ws.Cells.Item(420, 12).Value = 18000
$ws.Cells.Item(420, 13).Value = 18000
$ws.Cells.Item(420, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(420, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(420, 16).Value = 300
$ws.Cells.Item(420, 17).Value = 60
$ws.Cells.Item(420, 18).Value = "Hortaliza"
